# Refactor the lecturer and module version
# - Renumber "Semester" values in column A for rows 102-401
#   (rows 2-101 stay as 1; 102-201: 2->1; 202-301: 3->2; 302-401: 4->2)
# - Remove the old "Semester 5" block (rows 402-501)
# - Update the visible selection to match the new sheet state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the semester column for the affected ranges (top-down so we
# never overwrite a range we still need to read).
$ws.Range("A102:A201").Value = 1
$ws.Range("A202:A301").Value = 2
$ws.Range("A302:A401").Value = 2

# Drop the last 100 rows (previously "Semester 5" students) entirely,
# shifting everything below them up (there is nothing below, so this is
# simply a deletion of that block).
$ws.Range("A402:A501").EntireRow.Delete()

# Update the window selection to the new location referenced in the file.
$ws.Activate()
$ws.Range("G211").Select()
